# Fix a bug in FormalAmpersand Terms.xlsx:
#  1. The "Atoms" worksheet was an obsolete/duplicate sheet - remove it.
#  2. The "Rules" sheet was missing the blank separator row (row 3) that
#     every other example table in this workbook uses between the
#     attribute/type header rows and the example data row - insert it,
#     pushing the example row down from row 3 to row 4.
#  3. The "Terms" sheet had its last example (t4 / "r ISC s;t") stranded in
#     columns C:D instead of A:B like every other row in that table - move
#     it back into A6:B6.
#  4. The "Intersections" sheet had its example row stranded in columns
#     E:G instead of A:C (compare with the analogous "Compositions" sheet,
#     which correctly uses A:C) - move it back into A6:C6.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# 1. Remove the obsolete "Atoms" worksheet.
$wb.Worksheets("Atoms").Delete()

# 2. Rules sheet: insert a blank row above the example row (old row 3),
#    so the example data (r2 / Braga / t1 / t5) moves down to row 4.
$wsRules = $wb.Worksheets("Rules")
$wsRules.Range("A4:D4").Value2 = $wsRules.Range("A3:D3").Value2
$wsRules.Range("A3:D3").ClearContents()

# 3. Terms sheet: move the t4 row from C6:D6 back to A6:B6.
$wsTerms = $wb.Worksheets("Terms")
$wsTerms.Range("A6:B6").Value2 = $wsTerms.Range("C6:D6").Value2
$wsTerms.Range("C6:D6").ClearContents()

# 4. Intersections sheet: move the example row from E6:G6 back to A6:C6.
$wsIntersections = $wb.Worksheets("Intersections")
$wsIntersections.Range("A6:C6").Value2 = $wsIntersections.Range("E6:G6").Value2
$wsIntersections.Range("E6:G6").ClearContents()

# Restore/update each touched sheet's selection to where the edit landed.
$wsRules.Range("F6:G6").Select()
$wsTerms.Range("A6:B6").Select()
$wb.Worksheets("Compositions").Range("A6:C6").Select()
$wsIntersections.Range("A6:C6").Select()

# The workbook now has the "Rules" sheet active/selected.
$wsRules.Activate()
